# Add a new "Save" column (H) to the s_vals worksheet, matching the
# header formatting already used by the other header cells (G1 "sum", etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (G1, style index 1:
# bold, centered/top aligned, thin border) onto the new header cell H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the data value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
